$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 42641.540694444448
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 9952.56
$ws.Range("D4").Value = 9976.5
$ws.Range("E4").Value = 79.32
$ws.Range("F4").Value = 78.94
$ws.Range("G4").Value = $false
$ws.Range("G4").NumberFormat = "m/d/yy h:mm"
$ws.Range("H4").Value = -0.48
$ws.Range("I4").Value = $false
